# Final coding updates for the night
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codebook")

# --- Update column A width (25 -> 26) ---
# The runtime stores ColumnWidth in the XML "width" attribute as
# ColumnWidth + 0.8333333333333334 (standard Excel column-width padding),
# so compensate to land exactly on a stored width of 26.
$ws.Columns.Item(1).ColumnWidth = 26 - 0.8333333333333334

# --- Update I36 sample values ---
$ws.Range("I36").Value = "[2, 1, 5, 4]"

# --- Insert 13 new rows before row 70 to make room for the new CE definition
#     variables; this shifts the existing Media_/News_ rows (70-82) down to
#     rows 83-95, matching the shifted block in the target workbook. ---
$ws.Range("A70:A82").EntireRow.Insert()

# --- Populate the newly inserted rows (70-82) with the CE definition items ---
$newRows = @(
    @(70, "CEdef_voting",             "CE definition: Voting",                                                   "[1, 0]"),
    @(71, "CEdef_volunteering",       "CE definition: Volunteering or community service",                        "[0, 1]"),
    @(72, "CEdef_talking_community",  "CE definition: Talking with people in the community",                     "[0, 1]"),
    @(73, "CEdef_activism",           "CE definition: Activism or awareness raising",                            "[0, 1]"),
    @(74, "CEdef_political",          "CE definition: Political engagement",                                     "[0, 1]"),
    @(75, "CEdef_protesting",         "CE definition: Protesting or rallying",                                   "[0, 1]"),
    @(76, "CEdef_teaching",           "CE definition: Teaching or mentorship",                                   "[0, 1]"),
    @(77, "CEdef_fundraising",        "CE definition: Fundraising",                                              "[0, 1]"),
    @(78, "CEdef_arts_culture",       "CE definition: Attending community arts and culture events",              "[0, 1]"),
    @(79, "CEdef_community_business", "CE definition: Contributing to community businesses with a social cause", "[0, 1]"),
    @(80, "CEdef_donations",          "CE definition: Making donations to charities",                            "[0, 1]"),
    @(81, "CEdef_research",           "CE definition: Participating in research with community members",         "[0, 1]"),
    @(82, "CEdef_service_learning",   "CE definition: Service-learning experiences with a class",                "[0, 1]")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = "—"
    $ws.Cells.Item($r, 5).Value = "binary OHE"
    $ws.Cells.Item($r, 6).Value = "1=selected, 0=not selected, -9=missing"
    $ws.Cells.Item($r, 7).Value = -9
    $ws.Cells.Item($r, 8).Value = $row[2]
    $ws.Cells.Item($r, 9).Value = $row[3]
}
